# "Add files via upload" — re-saved workbook with:
#   1. The first sheet renamed from "Datos" to "Data".
#   2. The active cell on that sheet moved from I4 to F23
#      (the selection that gets persisted into the sheet's <sheetView>).

$wb = $excel.ActiveWorkbook

# 1) Rename the "Datos" sheet to "Data".
$ws = $wb.Worksheets.Item("Datos")
$ws.Name = "Data"

# 2) Move/persist the selection to F23 on that same (now renamed) sheet.
$ws.Activate()
$ws.Range("F23").Select()
